$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "401÷6=66, 5" "703÷5=140, 3"
Replace-Text "913÷6=152, 1" "219÷3=73, 0"
Replace-Text "409÷6=68, 1" "885÷5=177, 0"
Replace-Text "198÷4=49, 2" "109÷8=13, 5"
Replace-Text "939÷3=313, 0" "686÷3=228, 2"
Replace-Text "548÷5=109, 3" "941÷2=470, 1"
Replace-Text "274÷4=68, 2" "637÷7=91, 0"
Replace-Text "460÷9=51, 1" "122÷6=20, 2"
Replace-Text "664÷6=110, 4" "808÷7=115, 3"
Replace-Text "852÷2=426, 0" "900÷6=150, 0"
Replace-Text "483÷7=69, 0" "762÷9=84, 6"
Replace-Text "440÷5=88, 0" "499÷4=124, 3"
Replace-Text "561÷2=280, 1" "501÷7=71, 4"
Replace-Text "437÷9=48, 5" "309÷4=77, 1"
Replace-Text "636÷4=159, 0" "361÷8=45, 1"
Replace-Text "145÷9=16, 1" "847÷3=282, 1"
Replace-Text "382÷3=127, 1" "860÷9=95, 5"
Replace-Text "578÷7=82, 4" "398÷7=56, 6"
Replace-Text "255÷3=85, 0" "742÷4=185, 2"
Replace-Text "915÷6=152, 3" "869÷2=434, 1"
Replace-Text "209÷3=69, 2" "190÷6=31, 4"
Replace-Text "304÷7=43, 3" "228÷8=28, 4"
Replace-Text "185÷5=37, 0" "734÷2=367, 0"
Replace-Text "667÷9=74, 1" "591÷9=65, 6"
Replace-Text "806÷2=403, 0" "239÷8=29, 7"
